$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Update the feature-image prompt paragraph at the very end of the
# document (was the "Discover gameplay features..." meta description, now
# becomes the AI image generation prompt). Do this BEFORE inserting the new
# "Meta description" paragraph near the top so the Find/Replace below can't
# accidentally also match inside that newly-inserted paragraph.
# ---------------------------------------------------------------------------
$oldMetaText = "Discover gameplay features, pros and cons, and RTP in our review of Coils of Cash. Play for free and enter the free spins mode to increase your chances of winning big."
$newImagePrompt = "Create a feature image for Coils of Cash that captures the electrifying energy of the game. The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be wearing a tool belt filled with electrical equipment, showcasing the theme of the game. The background should be a vibrant blue, with electrical bolts and sparks surrounding the warrior. The warrior should also be surrounded by the game symbols, cascading down towards the bottom of the image. This feature image should exemplify the excitement and fun of playing Coils of Cash."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldMetaText, $true, $false, $false, $false, $false, $true, 1, $false, $newImagePrompt, 2)

# ---------------------------------------------------------------------------
# Step 2: Remove the duplicated bold "Play Coils of Cash Free: Review &
# Gameplay Features" paragraph that used to sit right before the paragraph
# handled above (it is now redundant with the title paragraph at the top,
# since a "Meta description" paragraph is being added there instead).
# ---------------------------------------------------------------------------
$dupTitleIndex = $d.Paragraphs.Count - 1
$dupTitlePara = $d.Paragraphs.Item($dupTitleIndex)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# Step 3: Insert a new "Meta description" paragraph right after the document
# title (Heading1) paragraph at the top.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleEnd = $titlePara.Range
$titleEnd.Collapse(0)            # wdCollapseEnd
$titleEnd.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Style = "Normal"

$boldLabel = "Meta description"
$restOfText = ": Discover gameplay features, pros and cons, and RTP in our review of Coils of Cash. Play for free and enter the free spins mode to increase your chances of winning big."

$insertionPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertionPoint.InsertAfter($boldLabel + $restOfText)

$boldRange = $d.Range($insertionPoint.Start, $insertionPoint.Start + $boldLabel.Length)
$boldRange.Bold = 1
